$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert the two new inventory rows.
#    - A new row at sheet row 10 ("multi-grip camber bar") which
#      shifts the old rows 10-30 down to 11-31.
#    - A new row that becomes row 31 ("foam roller"), inserted just
#      above the old blank/Total rows, shifting them down to 32/33.
# ------------------------------------------------------------------
$ws.Rows("10:10").Insert()
$ws.Rows("31:31").Insert()

# ------------------------------------------------------------------
# 2. Fill in the values for the two new rows. The shared-strings
#    table records new unique strings in the order they are first
#    written, so "foam roller" must be written before
#    "multi-grip camber bar" to land on the expected indices.
# ------------------------------------------------------------------
$ws.Range("A31").Value = "foam roller"
$ws.Range("B31").Value = 20

$ws.Range("A10").Value = "multi-grip camber bar"
$ws.Range("B10").Value = 194

# ------------------------------------------------------------------
# 3. Fix up the Total formula (now on row 33) to include the new rows.
# ------------------------------------------------------------------
$ws.Range("B33").Formula = "=SUM(B2:B31)"

# ------------------------------------------------------------------
# 4. Relocate the two cell comments so they stay attached to the same
#    logical items ("weights - iron" / "weights - bumper"), which
#    moved from B13/B14 down to B14/B15 because of the row insert at
#    row 10. Move the lower one first to avoid a collision.
# ------------------------------------------------------------------
$lowerText = $ws.Range("B14").Comment.Text()
$ws.Range("B14").Comment.Delete()
$ws.Range("B15").AddComment($lowerText)

$upperText = $ws.Range("B13").Comment.Text()
$ws.Range("B13").Comment.Delete()
$ws.Range("B14").AddComment($upperText)

# ------------------------------------------------------------------
# 5. Misc view-state changes captured in the workbook XML.
# ------------------------------------------------------------------
$ws.Range("D7").Select()

$excel.Width = 21435
$excel.Height = 11025
